# Applies odds updates to the FlashScore weekly games sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("J2").Value = 1.04
$ws.Range("L2").Value = 1.22
$ws.Range("N2").Value = 1.73

# Row 3
$ws.Range("G3").Value = 1.36
$ws.Range("U3").Value = 9

# Row 8
$ws.Range("I8").Value = 1.95
$ws.Range("J8").Value = 1.05
$ws.Range("K8").Value = 11
$ws.Range("N8").Value = 1.83
$ws.Range("O8").Value = 1.98

# Row 11
$ws.Range("G11").Value = 2.6
$ws.Range("H11").Value = 3.1
$ws.Range("I11").Value = 2.7
$ws.Range("L11").Value = 1.4
$ws.Range("M11").Value = 2.75
$ws.Range("N11").Value = 2.25
$ws.Range("O11").Value = 1.62
$ws.Range("T11").Value = 7.5
$ws.Range("U11").Value = 12
$ws.Range("V11").Value = 11
$ws.Range("W11").Value = 26
$ws.Range("X11").Value = 23
$ws.Range("AB11").Value = 15
$ws.Range("AE11").Value = 8
$ws.Range("AF11").Value = 13
$ws.Range("AG11").Value = 11
$ws.Range("AH11").Value = 29
$ws.Range("AI11").Value = 23
$ws.Range("AJ11").Value = 34
